$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(200).EntireRow.Insert()

$ws.Range("A200").Value = 10
$ws.Range("B200").Value = "Vega Modelo de Temuco"
$ws.Range("C200").Value = "La Araucanía"
$ws.Range("D200").Value = 44491
$ws.Range("E200").Value = 9
$ws.Range("F200").Value = 100112037
$ws.Range("G200").Value = "Cebollín"
$ws.Range("H200").Value = "Sin especificar"
$ws.Range("I200").Value = "Primera"
$ws.Range("J200").Value = 65
$ws.Range("K200").Value = 7000
$ws.Range("L200").Value = 7000
$ws.Range("M200").Value = 7000
$ws.Range("N200").Value = "$/docena de paquetes"
$ws.Range("O200").Value = "Provincia de Cautín"
$ws.Range("P200").Value = 583
$ws.Range("Q200").Value = 12
$ws.Range("R200").Value = "Hortaliza"
